$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (row 1) to the new lower-case naming convention.
# Order matters for shared-string table layout: percentile/gender columns
# (G:L) are updated first, then the summary-stat columns (B:F).
$ws.Range("G1").Value = "percent5th"
$ws.Range("H1").Value = "percent25th"
$ws.Range("I1").Value = "percent50th"
$ws.Range("J1").Value = "percent75th"
$ws.Range("K1").Value = "percent95th"
$ws.Range("L1").Value = "gender"
$ws.Range("B1").Value = "n"
$ws.Range("C1").Value = "min"
$ws.Range("D1").Value = "max"
$ws.Range("E1").Value = "mean"
$ws.Range("F1").Value = "sd"

# Update the selection / view to match the new workbook state.
$ws.Range("O4").Select()
